$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 20, shifting existing rows 20:111 down to 23:114
$ws.Rows("20:22").Insert()

# Populate the 3 newly inserted rows (20, 21, 22) with new data
# Row 20: Especial, bandeja 10 kilos
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44847
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107002
$ws.Range("J20").Value = "Chirimoya"
$ws.Range("K20").Value = "Cultivar IV Región"
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("Q20").Value = "$/bandeja 10 kilos"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 2050
$ws.Range("T20").Value = 10

# Row 21: Primera, bandeja 10 kilos
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44847
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = "Otros"
$ws.Range("I21").Value = 100107002
$ws.Range("J21").Value = "Chirimoya"
$ws.Range("K21").Value = "Cultivar IV Región"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 360
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 19000
$ws.Range("P21").Value = 18500
$ws.Range("Q21").Value = "$/bandeja 10 kilos"
$ws.Range("R21").Value = "Provincia de Limarí"
$ws.Range("S21").Value = 1850
$ws.Range("T21").Value = 10

# Row 22: Segunda, bandeja 10 kilos
$ws.Range("A22").Value = 2
$ws.Range("B22").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44847
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100107
$ws.Range("H22").Value = "Otros"
$ws.Range("I22").Value = 100107002
$ws.Range("J22").Value = "Chirimoya"
$ws.Range("K22").Value = "Cultivar IV Región"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 15500
$ws.Range("Q22").Value = "$/bandeja 10 kilos"
$ws.Range("R22").Value = "Provincia de Limarí"
$ws.Range("S22").Value = 1550
$ws.Range("T22").Value = 10

$ws.Range("D20:D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$addr = $ws.UsedRange.Address()
Write-Host "UsedRange: $addr"
